# Two new homogenization-weighting methods ("Holden" and "Rizzie Spiral") were
# added to the simulation right after "Spiral5", "Thomas Hex" was renamed to
# "Matthies Hex", and the whole averaging sweep (rows 2-31, columns C-T) was
# rerun, so every numeric result cell shifts/changes. Rebuild the entire data
# block (A2:T31) in one shot from the values produced by the rerun.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 30,20

# Row 2
$data[0,0] = 0
$data[0,1] = "HKL"
$data[0,2] = "[3, 2, 1]"
$data[0,3] = "[2, 2, 2]"
$data[0,4] = "[3, 1, 0]"
$data[0,5] = "[1, 1, 0]"
$data[0,6] = "[2, 2, 0]"
$data[0,7] = "[2, 0, 0]"
$data[0,8] = "[4, 0, 0]"
$data[0,9] = "[2, 1, 1]"
$data[0,10] = "1Pair-A"
$data[0,11] = "1Pair-B"
$data[0,12] = "2Pairs-A"
$data[0,13] = "2Pairs-B"
$data[0,14] = "3Pairs-A"
$data[0,15] = "3Pairs-B"
$data[0,16] = "3Pairs-C"
$data[0,17] = "4Pairs"
$data[0,18] = "5A4F"
$data[0,19] = "MaxUnique"

# Row 3
$data[1,0] = 1
$data[1,1] = "Spiral5"
$data[1,2] = 1.000033839275384
$data[1,3] = 1.000090237688475
$data[1,4] = 0.9999255538080086
$data[1,5] = 1.000033839275384
$data[1,6] = 1.000033839275384
$data[1,7] = 0.9998646451299626
$data[1,8] = 0.9998646451299626
$data[1,9] = 1.000033839275384
$data[1,10] = 1.000033839275384
$data[1,11] = 1.000033839275384
$data[1,12] = 0.9999492422026734
$data[1,13] = 0.9999492422026734
$data[1,14] = 0.9999413460711185
$data[1,15] = 0.9999774412269105
$data[1,16] = 0.9999774412269105
$data[1,17] = 0.999991540739029
$data[1,18] = 0.999991540739029
$data[1,19] = 0.9999969924087666

# Row 4
$data[2,0] = 2
$data[2,1] = "Holden"
$data[2,2] = 1.00163426422744
$data[2,3] = 1.004358027154505
$data[2,4] = 0.9964046266582384
$data[2,5] = 1.00163426422744
$data[2,6] = 1.00163426422744
$data[2,7] = 0.9934629589275957
$data[2,8] = 0.9934629589275957
$data[2,9] = 1.00163426422744
$data[2,10] = 1.00163426422744
$data[2,11] = 1.00163426422744
$data[2,12] = 0.9975486115775178
$data[2,13] = 0.9975486115775178
$data[2,14] = 0.9971672832710913
$data[2,15] = 0.9989104957941585
$data[2,16] = 0.9989104957941585
$data[2,17] = 0.9995914379024788
$data[2,18] = 0.9995914379024788
$data[2,19] = 0.9998547342371098

# Row 5
$data[3,0] = 3
$data[3,1] = "Rizzie Spiral"
$data[3,2] = 1.001126469056711
$data[3,3] = 1.003003900059562
$data[3,4] = 0.997521775932388
$data[3,5] = 1.001126469056711
$data[3,6] = 1.001126469056711
$data[3,7] = 0.9954941504626742
$data[3,8] = 0.9954941504626742
$data[3,9] = 1.001126469056711
$data[3,10] = 1.001126469056711
$data[3,11] = 1.001126469056711
$data[3,12] = 0.9983103097596926
$data[3,13] = 0.9983103097596926
$data[3,14] = 0.9980474651505911
$data[3,15] = 0.9992490295253654
$data[3,16] = 0.9992490295253654
$data[3,17] = 0.9997183894082018
$data[3,18] = 0.9997183894082018
$data[3,19] = 0.9998998722707929

# Row 6
$data[4,0] = 4
$data[4,1] = "RotRing OmegaMax-90"
$data[4,2] = 1.000612088006515
$data[4,3] = 1.001632231608403
$data[4,4] = 0.9986534106770185
$data[4,5] = 1.000612088006515
$data[4,6] = 1.000612088006515
$data[4,7] = 0.9975516571293739
$data[4,8] = 0.9975516571293739
$data[4,9] = 1.000612088006515
$data[4,10] = 1.000612088006515
$data[4,11] = 1.000612088006515
$data[4,12] = 0.9990818725679442
$data[4,13] = 0.9990818725679442
$data[4,14] = 0.9989390519376357
$data[4,15] = 0.9995919443808011
$data[4,16] = 0.9995919443808011
$data[4,17] = 0.9998469802872294
$data[4,18] = 0.9998469802872294
$data[4,19] = 0.9999455939057232

# Row 7
$data[5,0] = 5
$data[5,1] = "Equal Angle"
$data[5,2] = 1.00194845585015
$data[5,3] = 1.00519586911383
$data[5,4] = 0.9957134066570609
$data[5,5] = 1.00194845585015
$data[5,6] = 1.00194845585015
$data[5,7] = 0.9922061952665672
$data[5,8] = 0.9922061952665672
$data[5,9] = 1.00194845585015
$data[5,10] = 1.00194845585015
$data[5,11] = 1.00194845585015
$data[5,12] = 0.9970773255583587
$data[5,13] = 0.9970773255583587
$data[5,14] = 0.9966226859245927
$data[5,15] = 0.9987010356556225
$data[5,16] = 0.9987010356556225
$data[5,17] = 0.9995128907042544
$data[5,18] = 0.9995128907042544
$data[5,19] = 0.9998268064313182

# Row 8
$data[6,0] = 6
$data[6,1] = "Tilt Rotate"
$data[6,2] = 1.006536407967167
$data[6,3] = 1.017430349502477
$data[6,4] = 0.9856199511768708
$data[6,5] = 1.006536407967167
$data[6,6] = 1.006536407967167
$data[6,7] = 0.9738544553948659
$data[6,8] = 0.9738544553948659
$data[6,9] = 1.006536407967167
$data[6,10] = 1.006536407967167
$data[6,11] = 1.006536407967167
$data[6,12] = 0.9901954316810166
$data[6,13] = 0.9901954316810166
$data[6,14] = 0.988670271512968
$data[6,15] = 0.9956424237764002
$data[6,16] = 0.9956424237764002
$data[6,17] = 0.998365919824092
$data[6,18] = 0.998365919824092
$data[6,19] = 0.9994189966626194

# Row 9
$data[7,0] = 7
$data[7,1] = "CLR"
$data[7,2] = 1.00020332152893
$data[7,3] = 1.000542191497408
$data[7,4] = 0.9995526921759467
$data[7,5] = 1.00020332152893
$data[7,6] = 1.00020332152893
$data[7,7] = 0.9991867150210981
$data[7,8] = 0.9991867150210981
$data[7,9] = 1.00020332152893
$data[7,10] = 1.00020332152893
$data[7,11] = 1.00020332152893
$data[7,12] = 0.999695018275014
$data[7,13] = 0.999695018275014
$data[7,14] = 0.9996475762419915
$data[7,15] = 0.9998644526929859
$data[7,16] = 0.9998644526929859
$data[7,17] = 0.9999491699019718
$data[7,18] = 0.9999491699019718
$data[7,19] = 0.9999819272135402

# Row 10
$data[8,0] = 8
$data[8,1] = "Rizzie Hex"
$data[8,2] = 1.000016285896783
$data[8,3] = 1.00004342935799
$data[8,4] = 0.99996417135028
$data[8,5] = 1.000016285896783
$data[8,6] = 1.000016285896783
$data[8,7] = 0.999934858927486
$data[8,8] = 0.999934858927486
$data[8,9] = 1.000016285896783
$data[8,10] = 1.000016285896783
$data[8,11] = 1.000016285896783
$data[8,12] = 0.9999755724121346
$data[8,13] = 0.9999755724121346
$data[8,14] = 0.9999717720581831
$data[8,15] = 0.9999891435736842
$data[8,16] = 0.9999891435736842
$data[8,17] = 0.9999959291544589
$data[8,18] = 0.9999959291544589
$data[8,19] = 0.9999985528876842

# Row 11
$data[9,0] = 9
$data[9,1] = "Matthies Hex"
$data[9,2] = 1.000359495711253
$data[9,3] = 1.000958658147354
$data[9,4] = 0.9992091084844512
$data[9,5] = 1.000359495711253
$data[9,6] = 1.000359495711253
$data[9,7] = 0.9985620162466361
$data[9,8] = 0.9985620162466361
$data[9,9] = 1.000359495711253
$data[9,10] = 1.000359495711253
$data[9,11] = 1.000359495711253
$data[9,12] = 0.9994607559789446
$data[9,13] = 0.9994607559789446
$data[9,14] = 0.9993768734807801
$data[9,15] = 0.9997603358897141
$data[9,16] = 0.9997603358897141
$data[9,17] = 0.9999101258450989
$data[9,18] = 0.9999101258450989
$data[9,19] = 0.9999680450020335

# Row 12
$data[10,0] = 10
$data[10,1] = "Tilt Rotate_Partial"
$data[10,2] = 1.006633693420423
$data[10,3] = 1.017689774119754
$data[10,4] = 0.9854059253894013
$data[10,5] = 1.006633693420423
$data[10,6] = 1.006633693420423
$data[10,7] = 0.9734653173481661
$data[10,8] = 0.9734653173481661
$data[10,9] = 1.006633693420423
$data[10,10] = 1.006633693420423
$data[10,11] = 1.006633693420423
$data[10,12] = 0.9900495053842946
$data[10,13] = 0.9900495053842946
$data[10,14] = 0.9885016453859968
$data[10,15] = 0.995577568063004
$data[10,16] = 0.995577568063004
$data[10,17] = 0.9983415994023588
$data[10,18] = 0.9983415994023588
$data[10,19] = 0.9994103495197652

# Row 13
$data[11,0] = 11
$data[11,1] = "RotRing OmegaMax-60"
$data[11,2] = 1.000577015440336
$data[11,3] = 1.001538707125926
$data[11,4] = 0.9987305681677905
$data[11,5] = 1.000577015440336
$data[11,6] = 1.000577015440336
$data[11,7] = 0.9976919445798013
$data[11,8] = 0.9976919445798013
$data[11,9] = 1.000577015440336
$data[11,10] = 1.000577015440336
$data[11,11] = 1.000577015440336
$data[11,12] = 0.9991344800100685
$data[11,13] = 0.9991344800100685
$data[11,14] = 0.9989998427293093
$data[11,15] = 0.999615325153491
$data[11,16] = 0.999615325153491
$data[11,17] = 0.9998557477252022
$data[11,18] = 0.9998557477252022
$data[11,19] = 0.9999487110324211

# Row 14
$data[12,0] = 12
$data[12,1] = "Equal Angle_Partial"
$data[12,2] = 1.001990686431583
$data[12,3] = 1.005308476694737
$data[12,4] = 0.9956205035999998
$data[12,5] = 1.001990686431583
$data[12,6] = 1.001990686431583
$data[12,7] = 0.9920372804315766
$data[12,8] = 0.9920372804315766
$data[12,9] = 1.001990686431583
$data[12,10] = 1.001990686431583
$data[12,11] = 1.001990686431583
$data[12,12] = 0.9970139834315797
$data[12,13] = 0.9970139834315797
$data[12,14] = 0.9965494901543863
$data[12,15] = 0.9986728844315808
$data[12,16] = 0.9986728844315808
$data[12,17] = 0.9995023349315812
$data[12,18] = 0.9995023349315812
$data[12,19] = 0.9998230533368435

# Row 15
$data[13,0] = 13
$data[13,1] = "Rizzie Hex_Partial"
$data[13,2] = 0.998695330822546
$data[13,3] = 0.9965208801138797
$data[13,4] = 1.002870273343795
$data[13,5] = 0.998695330822546
$data[13,6] = 0.998695330822546
$data[13,7] = 1.005218681497794
$data[13,8] = 1.005218681497794
$data[13,9] = 0.998695330822546
$data[13,10] = 0.998695330822546
$data[13,11] = 0.998695330822546
$data[13,12] = 1.00195700616017
$data[13,13] = 1.00195700616017
$data[13,14] = 1.002261428554712
$data[13,15] = 1.000869781047629
$data[13,16] = 1.000869781047629
$data[13,17] = 1.000326168491358
$data[13,18] = 1.000326168491358
$data[13,19] = 1.000115971237185

# Row 16
$data[14,0] = 14
$data[14,1] = "ND Single"
$data[14,2] = 1.011461099999999
$data[14,3] = 1.030562799999999
$data[14,4] = 0.9747856700000013
$data[14,5] = 1.011461099999999
$data[14,6] = 1.011461099999999
$data[14,7] = 0.9541557599999997
$data[14,8] = 0.9541557599999997
$data[14,9] = 1.011461099999999
$data[14,10] = 1.011461099999999
$data[14,11] = 1.011461099999999
$data[14,12] = 0.9828084299999995
$data[14,13] = 0.9828084299999995
$data[14,14] = 0.9801341766666667
$data[14,15] = 0.9923593199999994
$data[14,16] = 0.9923593199999994
$data[14,17] = 0.9971347649999993
$data[14,18] = 0.9971347649999993
$data[14,19] = 0.9989812549999996

# Row 17
$data[15,0] = 15
$data[15,1] = "RD Single"
$data[15,2] = 1.0042979
$data[15,3] = 1.0114611
$data[15,4] = 0.99054463
$data[15,5] = 1.0042979
$data[15,6] = 1.0042979
$data[15,7] = 0.98280841
$data[15,8] = 0.98280841
$data[15,9] = 1.0042979
$data[15,10] = 1.0042979
$data[15,11] = 1.0042979
$data[15,12] = 0.9935531550000001
$data[15,13] = 0.9935531550000001
$data[15,14] = 0.9925503133333334
$data[15,15] = 0.9971347366666667
$data[15,16] = 0.9971347366666667
$data[15,17] = 0.9989255275000001
$data[15,18] = 0.9989255275000001
$data[15,19] = 0.9996179733333334

# Row 18
$data[16,0] = 16
$data[16,1] = "TD Single"
$data[16,2] = 1.0042976
$data[16,3] = 1.0114603
$data[16,4] = 0.99054523
$data[16,5] = 1.0042976
$data[16,6] = 1.0042976
$data[16,7] = 0.98280951
$data[16,8] = 0.98280951
$data[16,9] = 1.0042976
$data[16,10] = 1.0042976
$data[16,11] = 1.0042976
$data[16,12] = 0.993553555
$data[16,13] = 0.993553555
$data[16,14] = 0.9925507800000001
$data[16,15] = 0.9971349033333333
$data[16,16] = 0.9971349033333333
$data[16,17] = 0.9989255774999999
$data[16,18] = 0.9989255774999999
$data[16,19] = 0.9996179733333332

# Row 19
$data[17,0] = 17
$data[17,1] = "Morris Single"
$data[17,2] = 0.9966870399999999
$data[17,3] = 0.9911654300000001
$data[17,4] = 1.0072885
$data[17,5] = 0.9966870399999999
$data[17,6] = 0.9966870399999999
$data[17,7] = 1.0132519
$data[17,8] = 1.0132519
$data[17,9] = 0.9966870399999999
$data[17,10] = 0.9966870399999999
$data[17,11] = 0.9966870399999999
$data[17,12] = 1.00496947
$data[17,13] = 1.00496947
$data[17,14] = 1.00574248
$data[17,15] = 1.00220866
$data[17,16] = 1.00220866
$data[17,17] = 1.000828255
$data[17,18] = 1.000828255
$data[17,19] = 1.000294491666667

# Row 20
$data[18,0] = 18
$data[18,1] = "Ring Perpendicular to ND"
$data[18,2] = 1.004297897260274
$data[18,3] = 1.011461067123287
$data[18,4] = 0.9905446256164382
$data[18,5] = 1.004297897260274
$data[18,6] = 1.004297897260274
$data[18,7] = 0.9828084100000002
$data[18,8] = 0.9828084100000002
$data[18,9] = 1.004297897260274
$data[18,10] = 1.004297897260274
$data[18,11] = 1.004297897260274
$data[18,12] = 0.993553153630137
$data[18,13] = 0.993553153630137
$data[18,14] = 0.992550310958904
$data[18,15] = 0.9971347348401826
$data[18,16] = 0.9971347348401826
$data[18,17] = 0.9989255254452054
$data[18,18] = 0.9989255254452054
$data[18,19] = 0.9996179657534245

# Row 21
$data[19,0] = 19
$data[19,1] = "Ring Perpendicular to RD"
$data[19,2] = 1.001941594736842
$data[19,3] = 1.005177570526316
$data[19,4] = 0.9957285084210525
$data[19,5] = 1.001941594736842
$data[19,6] = 1.001941594736842
$data[19,7] = 0.992233654736842
$data[19,8] = 0.992233654736842
$data[19,9] = 1.001941594736842
$data[19,10] = 1.001941594736842
$data[19,11] = 1.001941594736842
$data[19,12] = 0.9970876247368421
$data[19,13] = 0.9970876247368421
$data[19,14] = 0.9966345859649123
$data[19,15] = 0.998705614736842
$data[19,16] = 0.998705614736842
$data[19,17] = 0.9995146097368421
$data[19,18] = 0.9995146097368421
$data[19,19] = 0.9998274196491228

# Row 22
$data[20,0] = 20
$data[20,1] = "Ring Perpendicular to TD"
$data[20,2] = 1.001941502105263
$data[20,3] = 1.005177319473684
$data[20,4] = 0.9957287047368423
$data[20,5] = 1.001941502105263
$data[20,6] = 1.001941502105263
$data[20,7] = 0.9922340194736842
$data[20,8] = 0.9922340194736842
$data[20,9] = 1.001941502105263
$data[20,10] = 1.001941502105263
$data[20,11] = 1.001941502105263
$data[20,12] = 0.9970877607894735
$data[20,13] = 0.9970877607894735
$data[20,14] = 0.9966347421052631
$data[20,15] = 0.9987056745614034
$data[20,16] = 0.9987056745614034
$data[20,17] = 0.9995146314473682
$data[20,18] = 0.9995146314473682
$data[20,19] = 0.9998274249999999

# Row 23
$data[21,0] = 21
$data[21,1] = "OffsetFTD"
$data[21,2] = 0.9984981641398782
$data[21,3] = 0.9959951328860016
$data[21,4] = 1.003304018664318
$data[21,5] = 0.9984981641398782
$data[21,6] = 0.9984981641398782
$data[21,7] = 1.006007306339125
$data[21,8] = 1.006007306339125
$data[21,9] = 0.9984981641398782
$data[21,10] = 0.9984981641398782
$data[21,11] = 0.9984981641398782
$data[21,12] = 1.002252735239502
$data[21,13] = 1.002252735239502
$data[21,14] = 1.002603163047774
$data[21,15] = 1.001001211539627
$data[21,16] = 1.001001211539627
$data[21,17] = 1.00037544968969
$data[21,18] = 1.00037544968969
$data[21,19] = 1.00013349171818

# Row 24
$data[22,0] = 22
$data[22,1] = "OffsetATD"
$data[22,2] = 0.9996122454127948
$data[22,3] = 0.9989659978780547
$data[22,4] = 1.00085305578587
$data[22,5] = 0.9996122454127948
$data[22,6] = 0.9996122454127948
$data[22,7] = 1.001551005564818
$data[22,8] = 1.001551005564818
$data[22,9] = 0.9996122454127948
$data[22,10] = 0.9996122454127948
$data[22,11] = 0.9996122454127948
$data[22,12] = 1.000581625488806
$data[22,13] = 1.000581625488806
$data[22,14] = 1.000672102254494
$data[22,15] = 1.000258498796803
$data[22,16] = 1.000258498796803
$data[22,17] = 1.000096935450801
$data[22,18] = 1.000096935450801
$data[22,19] = 1.000034465911188

# Row 25
$data[23,0] = 23
$data[23,1] = "OffsetF45"
$data[23,2] = 0.9984980665675514
$data[23,3] = 0.995994846893131
$data[23,4] = 1.003304252380399
$data[23,5] = 0.9984980665675514
$data[23,6] = 0.9984980665675514
$data[23,7] = 1.006007732561564
$data[23,8] = 1.006007732561564
$data[23,9] = 0.9984980665675514
$data[23,10] = 0.9984980665675514
$data[23,11] = 0.9984980665675514
$data[23,12] = 1.002252899564557
$data[23,13] = 1.002252899564557
$data[23,14] = 1.002603350503171
$data[23,15] = 1.001001288565555
$data[23,16] = 1.001001288565555
$data[23,17] = 1.000375483066054
$data[23,18] = 1.000375483066054
$data[23,19] = 1.000133505256291

# Row 26
$data[24,0] = 24
$data[24,1] = "OffsetA45"
$data[24,2] = 0.9996122236349583
$data[24,3] = 0.9989659182172428
$data[24,4] = 1.000853110027781
$data[24,5] = 0.9996122236349583
$data[24,6] = 0.9996122236349583
$data[24,7] = 1.001551118814635
$data[24,8] = 1.001551118814635
$data[24,9] = 0.9996122236349583
$data[24,10] = 0.9996122236349583
$data[24,11] = 0.9996122236349583
$data[24,12] = 1.000581671224797
$data[24,13] = 1.000581671224797
$data[24,14] = 1.000672150825791
$data[24,15] = 1.000258522028184
$data[24,16] = 1.000258522028184
$data[24,17] = 1.000096947429878
$data[24,18] = 1.000096947429878
$data[24,19] = 1.000034469660756

# Row 27
$data[25,0] = 25
$data[25,1] = "OffsetFRD"
$data[25,2] = 0.9984979121965397
$data[25,3] = 0.9959944294882275
$data[25,4] = 1.003304584185431
$data[25,5] = 0.9984979121965397
$data[25,6] = 0.9984979121965397
$data[25,7] = 1.006008349699427
$data[25,8] = 1.006008349699427
$data[25,9] = 0.9984979121965397
$data[25,10] = 0.9984979121965397
$data[25,11] = 0.9984979121965397
$data[25,12] = 1.002253130947983
$data[25,13] = 1.002253130947983
$data[25,14] = 1.002603615360466
$data[25,15] = 1.001001391364169
$data[25,16] = 1.001001391364169
$data[25,17] = 1.000375521572262
$data[25,18] = 1.000375521572262
$data[25,19] = 1.000133516660451

# Row 28
$data[26,0] = 26
$data[26,1] = "OffsetARD"
$data[26,2] = 0.9996121799737203
$data[26,3] = 0.9989658086527516
$data[26,4] = 1.000853207938217
$data[26,5] = 0.9996121799737203
$data[26,6] = 0.9996121799737203
$data[26,7] = 1.001551287585155
$data[26,8] = 1.001551287585155
$data[26,9] = 0.9996121799737203
$data[26,10] = 0.9996121799737203
$data[26,11] = 0.9996121799737203
$data[26,12] = 1.000581733779438
$data[26,13] = 1.000581733779438
$data[26,14] = 1.000672225165697
$data[26,15] = 1.000258549177532
$data[26,16] = 1.000258549177532
$data[26,17] = 1.000096956876579
$data[26,18] = 1.000096956876579
$data[26,19] = 1.000034474016214

# Row 29
$data[27,0] = 27
$data[27,1] = "Gaussian Quadrature"
$data[27,2] = 1.000683084663636
$data[27,3] = 1.001821539769471
$data[27,4] = 0.9984972263116271
$data[27,5] = 1.000683084663636
$data[27,6] = 1.000683084663636
$data[27,7] = 0.9972676847071752
$data[27,8] = 0.9972676847071752
$data[27,9] = 1.000683084663636
$data[27,10] = 1.000683084663636
$data[27,11] = 1.000683084663636
$data[27,12] = 0.9989753846854055
$data[27,13] = 0.9989753846854055
$data[27,14] = 0.9988159985608127
$data[27,15] = 0.9995446180114822
$data[27,16] = 0.9995446180114822
$data[27,17] = 0.9998292346745206
$data[27,18] = 0.9998292346745206
$data[27,19] = 0.9999392841298634

# Row 30
$data[28,0] = 28
$data[28,1] = "Michael-CCHex"
$data[28,2] = 0.9998815122898956
$data[28,3] = 0.9996840048444549
$data[28,4] = 1.000260699731463
$data[28,5] = 0.9998815122898956
$data[28,6] = 0.9998815122898956
$data[28,7] = 1.000473996338852
$data[28,8] = 1.000473996338852
$data[28,9] = 0.9998815122898956
$data[28,10] = 0.9998815122898956
$data[28,11] = 0.9998815122898956
$data[28,12] = 1.000177754314374
$data[28,13] = 1.000177754314374
$data[28,14] = 1.000205402786737
$data[28,15] = 1.000079006972881
$data[28,16] = 1.000079006972881
$data[28,17] = 1.000029633302135
$data[28,18] = 1.000029633302135
$data[28,19] = 1.000010539630743

# Row 31
$data[29,0] = 29
$data[29,1] = "Michael-SNHex"
$data[29,2] = 0.9977960775245921
$data[29,3] = 0.9941228420082292
$data[29,4] = 1.00484865147597
$data[29,5] = 0.9977960775245921
$data[29,6] = 0.9977960775245921
$data[29,7] = 1.008815735464564
$data[29,8] = 1.008815735464564
$data[29,9] = 0.9977960775245921
$data[29,10] = 0.9977960775245921
$data[29,11] = 0.9977960775245921
$data[29,12] = 1.003305906494578
$data[29,13] = 1.003305906494578
$data[29,14] = 1.003820154821709
$data[29,15] = 1.001469296837916
$data[29,16] = 1.001469296837916
$data[29,17] = 1.000550992009585
$data[29,18] = 1.000550992009585
$data[29,19] = 1.000195910253757

$ws.Range("A2:T31").Value = $data

# The two brand-new rows need the same bold/centered/bordered format that
# column A already uses elsewhere in the table; clone it from the row above.
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)
$excel.CutCopyMode = $false
